$wb = $excel.ActiveWorkbook

# Sheet1: "targets" - columns A (id), B (name)
$ws1 = $wb.Worksheets.Item("targets")
$ws1.Range("A5").Value = 782
$ws1.Range("B5").Value = "Chronic Thromboembolic Pulmonary Hypertension"

# Sheet2: "outcomes" - columns A (id), B (name), C (clean_window)
$ws2 = $wb.Worksheets.Item("outcomes")
$ws2.Range("A5").Value = 782
$ws2.Range("B5").Value = "Chronic Thromboembolic Pulmonary Hypertension"
$ws2.Range("C5").Value = 9999
